# Auto-generated script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '74.528.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +8.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.591.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '585.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.73%  '

$ws.Range("E7").Value = '  -0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.537'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.205'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +22.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.588.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.50%  '

$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.362'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.74%  '

$ws.Range("E13").Value = '  +4.79%  '

$ws.Range("E14").Value = '  +8.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.462.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.066.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.60%  '

$ws.Range("E17").Value = '  +13.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.598.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +32.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +12.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +11.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.97%  '

$ws.Range("E23").Value = '  +7.50%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.12%  '

$ws.Range("E26").Value = '  +13.62%  '

$ws.Range("E27").Value = '  +14.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.726.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.996'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0950'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +16.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +21.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.82%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '508.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("E36").Value = '  +14.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.24'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.89%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  +13.75%  '

$ws.Range("E42").Value = '  +12.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.327'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '156.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +18.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +9.12%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0856'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +19.59%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '38.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.03%  '

$ws.Range("E49").Value = '  +8.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.524'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +20.71%  '
